$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tasks to append (Arquitectura related use cases)
$tasks = @(
    "Gestión de Perfiles de Usuario",
    "Gestión de Login/Logout",
    "Gestión de Múltiples Idiomas",
    "Gestión de Bitácora y Control de Cambios",
    "Gestión de Backup",
    "Gestión de DV",
    "Gestión de Encriptado",
    "Casos de Prueba",
    "Manual de Instalación",
    "Ayuda en Línea"
)

$startRow = 15
for ($i = 0; $i -lt $tasks.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $tasks[$i]
}

$endRow = $startRow + $tasks.Length - 1

# Apply a thin box border + vertical-center alignment to the first new
# row, then propagate that exact formatting to the remaining new rows
# via a single format-only paste so the style table stays minimal
# (mirrors what Excel itself produces when the formatting is applied
# once and then copied down).
$firstCell = $ws.Range("A$startRow")
$firstCell.Borders.Weight = 2
$firstCell.VerticalAlignment = -4108

$firstCell.Copy()
$restRange = $ws.Range("A" + ($startRow + 1) + ":A$endRow")
$restRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A$startRow").Select()
